# Adds a new "2022-Q3" sheet (fund-holding detail) right after the "总计"
# summary sheet, populates it with the quarter's fund data, and updates the
# "总计" sheet with the new 2022-Q3 summary row (shifting the other quarters
# down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet immediately after "总计" (i.e. before
#    the worksheet that is currently in slot 2, "2022-Q2").
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Match the page-margin convention used by the sibling quarter sheets
# (left/right 0.75in, top/bottom 1in, header/footer 0.5in -> specified in
# points for the PageSetup API: 1in = 72pt).
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Header row.
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Fund-holding detail rows for 2022-Q3.
# Columns: index, code, name, scale, stockPosition, positionShare, marketValue, rank
$q3Rows = @(
    @(0,  "005314", "万家中证1000指数增强C",             "14.28", "94.11", "1.14", "0.1628", 1),
    @(1,  "005313", "万家中证1000指数增强A",             "13.25", "94.11", "1.14", "0.1510", 1),
    @(2,  "005457", "景顺长城量化小盘股票",               "6.57",  "93.58", "1.78", "0.1169", 3),
    @(3,  "006165", "建信中证1000指数增强A",             "3.87",  "84.02", "1.62", "0.0627", 2),
    @(4,  "014202", "天弘中证1000指数增强C",             "3.69",  "94.06", "1.58", "0.0583", 7),
    @(5,  "014201", "天弘中证1000指数增强A",             "3.68",  "94.06", "1.58", "0.0581", 7),
    @(6,  "015496", "景顺中证1000指数增强C",             "1.83",  "92.63", "1.69", "0.0309", 6),
    @(7,  "006166", "建信中证1000指数增强C",             "1.89",  "84.02", "1.62", "0.0306", 2),
    @(8,  "015495", "景顺中证1000指数增强A",             "0.69",  "92.63", "1.69", "0.0117", 6),
    @(9,  "015148", "华安中证1000指数增强A",             "1.50",  "91.03", "0.61", "0.0092", 10),
    @(10, "015149", "华安中证1000指数增强C",             "0.77",  "91.03", "0.61", "0.0047", 10),
    @(11, "005120", "上投摩根量化多因子灵活配置混合",       "0.19",  "92.91", "1.74", "0.0033", 10),
    @(12, "005167", "嘉实润泽量化一年定期开放混合",         "0.55",  "24.55", "0.58", "0.0032", 3),
    @(13, "013442", "建信中证1000指数增强E",             "0.18",  "84.02", "1.62", "0.0029", 2)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    # Prefix numeric-looking text (fund code / scale / position figures) with
    # an apostrophe so Excel stores it as text, exactly like the other
    # quarter sheets (fund codes such as "005314" must keep their leading
    # zero, and the percentage-like figures are text too).
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = "'" + $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add the 2022-Q3 row at the top of the
#    data (row 2) and shift the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q3", 14, 0.71),
    @(1, "2022-Q2", 7,  0.21),
    @(2, "2021-Q3", 1,  0.02),
    @(3, "2021-Q2", 1,  0.02),
    @(4, "2021-Q1", 1,  0.02),
    @(5, "2020-Q4", 3,  0.13)
)

$r = 2
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# The new last row (A7) falls outside the sheet's original used range, so it
# doesn't inherit the bordered/centered style ("s=2") the other index cells
# in column A use. Copy that formatting across from the row above it.
$summary.Cells.Item(6, 1).Copy()
$summary.Cells.Item(7, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Restore the original tab selection (the last sheet, "2020-Q4", was the
#    active tab before this edit).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
